$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the values of columns A, Q, R, AC between row 4 and row 5.

$cols = @("A", "Q", "R", "AC")

foreach ($col in $cols) {
    $addr4 = "$col`4"
    $addr5 = "$col`5"
    $val4 = $ws.Range($addr4).Value2
    $val5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value2 = $val5
    $ws.Range($addr5).Value2 = $val4
}
